# Auto-generated edit script: refresh market-data columns (H-N) across all Leve-profit sheets
# matching the scheduled runner's data update.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 2220.0476   # H28: 2219.9048 -> 2220.0476
$ws.Cells.Item(28, 9).Value = 751.8461   # I28: 751.61536 -> 751.8461
$ws.Cells.Item(28, 11).Value = 751.8461   # K28: 751.61536 -> 751.8461
$ws.Cells.Item(28, 13).Value = -266.8461   # M28: -266.61536 -> -266.8461
$ws.Cells.Item(33, 8).Value = 157.75   # H33: 170.28572 -> 157.75
$ws.Cells.Item(33, 9).Value = 157.75   # I33: 170.28572 -> 157.75
$ws.Cells.Item(33, 11).Value = 157.75   # K33: 170.28572 -> 157.75
$ws.Cells.Item(33, 13).Value = 71.25   # M33: 58.71428 -> 71.25
$ws.Cells.Item(51, 8).Value = 3187.5   # H51: 2458.5 -> 3187.5
$ws.Cells.Item(51, 9).Value = 2375   # I51: 1916.6666 -> 2375
$ws.Cells.Item(51, 10).Value = 4000   # J51: 3000.3333 -> 4000
$ws.Cells.Item(51, 11).Value = 2375   # K51: 1916.6666 -> 2375
$ws.Cells.Item(51, 12).Value = 4000   # L51: 3000.3333 -> 4000
$ws.Cells.Item(51, 13).Value = -1891   # M51: -1432.6666 -> -1891
$ws.Cells.Item(51, 14).Value = -4968   # N51: -3968.3333 -> -4968
$ws.Cells.Item(98, 8).Value = 680.34784   # H98: 708.2857 -> 680.34784
$ws.Cells.Item(98, 9).Value = 688.5909   # I98: 718.75 -> 688.5909
$ws.Cells.Item(98, 11).Value = 688.5909   # K98: 718.75 -> 688.5909
$ws.Cells.Item(98, 13).Value = 809.4091   # M98: 779.25 -> 809.4091
$ws.Cells.Item(111, 8).Value = 653.7273   # H111: 669.5 -> 653.7273
$ws.Cells.Item(111, 9).Value = 406.83334   # I111: 407.5 -> 406.83334
$ws.Cells.Item(111, 10).Value = 950   # J111: 1062.5 -> 950
$ws.Cells.Item(111, 11).Value = 1220.50002   # K111: 1222.5 -> 1220.50002
$ws.Cells.Item(111, 12).Value = 2850   # L111: 3187.5 -> 2850
$ws.Cells.Item(111, 13).Value = 1846.49998   # M111: 1844.5 -> 1846.49998
$ws.Cells.Item(111, 14).Value = -8984   # N111: -9321.5 -> -8984
$ws.Cells.Item(113, 8).Value = 4078.1052   # H113: 4165.778 -> 4078.1052
$ws.Cells.Item(113, 9).Value = 3783.6924   # I113: 3890.6667 -> 3783.6924
$ws.Cells.Item(113, 11).Value = 3783.6924   # K113: 3890.6667 -> 3783.6924
$ws.Cells.Item(113, 13).Value = -529.6923999999999   # M113: -636.6667000000002 -> -529.6923999999999
$ws.Cells.Item(122, 8).Value = 680.34784   # H122: 708.2857 -> 680.34784
$ws.Cells.Item(122, 9).Value = 688.5909   # I122: 718.75 -> 688.5909
$ws.Cells.Item(122, 11).Value = 2065.7727   # K122: 2156.25 -> 2065.7727
$ws.Cells.Item(122, 13).Value = 384.2273   # M122: 293.75 -> 384.2273
$ws.Cells.Item(125, 8).Value = 4517.2354   # H125: 4712.375 -> 4517.2354
$ws.Cells.Item(125, 9).Value = 2984.1428   # I125: 3248.1667 -> 2984.1428
$ws.Cells.Item(125, 10).Value = 5590.4   # J125: 5590.9 -> 5590.4
$ws.Cells.Item(125, 11).Value = 26857.2852   # K125: 29233.5003 -> 26857.2852
$ws.Cells.Item(125, 12).Value = 50313.6   # L125: 50318.1 -> 50313.6
$ws.Cells.Item(125, 13).Value = -24397.2852   # M125: -26773.5003 -> -24397.2852
$ws.Cells.Item(125, 14).Value = -55233.6   # N125: -55238.1 -> -55233.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 11080.206   # H32: 4827.8906 -> 11080.206
$ws.Cells.Item(32, 9).Value = 4311.476   # I32: 1673.7368 -> 4311.476
$ws.Cells.Item(32, 10).Value = 22014.309   # J32: 30511.715 -> 22014.309
$ws.Cells.Item(32, 11).Value = 4311.476   # K32: 1673.7368 -> 4311.476
$ws.Cells.Item(32, 12).Value = 22014.309   # L32: 30511.715 -> 22014.309
$ws.Cells.Item(32, 13).Value = -4024.476   # M32: -1386.7368 -> -4024.476
$ws.Cells.Item(32, 14).Value = -22588.309   # N32: -31085.715 -> -22588.309
$ws.Cells.Item(74, 8).Value = 35406.766   # H74: 35425.7 -> 35406.766
$ws.Cells.Item(74, 9).Value = 43120.957   # I74: 41450.84 -> 43120.957
$ws.Cells.Item(74, 10).Value = 4550   # J74: 5300 -> 4550
$ws.Cells.Item(74, 11).Value = 43120.957   # K74: 41450.84 -> 43120.957
$ws.Cells.Item(74, 12).Value = 4550   # L74: 5300 -> 4550
$ws.Cells.Item(74, 13).Value = -42246.957   # M74: -40576.84 -> -42246.957
$ws.Cells.Item(74, 14).Value = -6298   # N74: -7048 -> -6298
$ws.Cells.Item(77, 8).Value = 35406.766   # H77: 35425.7 -> 35406.766
$ws.Cells.Item(77, 9).Value = 43120.957   # I77: 41450.84 -> 43120.957
$ws.Cells.Item(77, 10).Value = 4550   # J77: 5300 -> 4550
$ws.Cells.Item(77, 11).Value = 215604.785   # K77: 207254.2 -> 215604.785
$ws.Cells.Item(77, 12).Value = 22750   # L77: 26500 -> 22750
$ws.Cells.Item(77, 13).Value = -211236.785   # M77: -202886.2 -> -211236.785
$ws.Cells.Item(77, 14).Value = -31486   # N77: -35236 -> -31486
$ws.Cells.Item(102, 8).Value = 71390.94   # H102: 51988.59 -> 71390.94
$ws.Cells.Item(102, 9).Value = 92296.09   # I102: 59808.883 -> 92296.09
$ws.Cells.Item(102, 11).Value = 92296.09   # K102: 59808.883 -> 92296.09
$ws.Cells.Item(102, 13).Value = -90674.09   # M102: -58186.883 -> -90674.09
$ws.Cells.Item(122, 8).Value = 3086.75   # H122: 2235.4 -> 3086.75
$ws.Cells.Item(122, 9).Value = 1798.5   # I122: 1700.5 -> 1798.5
$ws.Cells.Item(122, 11).Value = 5395.5   # K122: 5101.5 -> 5395.5
$ws.Cells.Item(122, 13).Value = -2945.5   # M122: -2651.5 -> -2945.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(8, 8).Value = 7000   # H8: 6000 -> 7000
$ws.Cells.Item(8, 9).Value = 7000   # I8: 6000 -> 7000
$ws.Cells.Item(8, 11).Value = 7000   # K8: 6000 -> 7000
$ws.Cells.Item(8, 13).Value = -6860   # M8: -5860 -> -6860
$ws.Cells.Item(11, 8).Value = 0   # H11: 14000 -> 0
$ws.Cells.Item(11, 10).Value = 0   # J11: 14000 -> 0
$ws.Cells.Item(11, 12).Value = 0   # L11: 14000 -> 0
$ws.Cells.Item(86, 8).Value = 3684.6956   # H86: 3675.15 -> 3684.6956
$ws.Cells.Item(86, 9).Value = 3947.1428   # I86: 3743.9285 -> 3947.1428
$ws.Cells.Item(86, 10).Value = 3276.4443   # J86: 3514.6667 -> 3276.4443
$ws.Cells.Item(86, 11).Value = 3947.1428   # K86: 3743.9285 -> 3947.1428
$ws.Cells.Item(86, 12).Value = 3276.4443   # L86: 3514.6667 -> 3276.4443
$ws.Cells.Item(86, 13).Value = -2824.1428   # M86: -2620.9285 -> -2824.1428
$ws.Cells.Item(86, 14).Value = -5522.4443   # N86: -5760.6667 -> -5522.4443
$ws.Cells.Item(89, 8).Value = 3684.6956   # H89: 3675.15 -> 3684.6956
$ws.Cells.Item(89, 9).Value = 3947.1428   # I89: 3743.9285 -> 3947.1428
$ws.Cells.Item(89, 10).Value = 3276.4443   # J89: 3514.6667 -> 3276.4443
$ws.Cells.Item(89, 11).Value = 19735.714   # K89: 18719.6425 -> 19735.714
$ws.Cells.Item(89, 12).Value = 16382.2215   # L89: 17573.3335 -> 16382.2215
$ws.Cells.Item(89, 13).Value = -14119.714   # M89: -13103.6425 -> -14119.714
$ws.Cells.Item(89, 14).Value = -27614.2215   # N89: -28805.3335 -> -27614.2215
$ws.Cells.Item(99, 8).Value = 2349131.8   # H99: 2487269.5 -> 2349131.8
$ws.Cells.Item(99, 9).Value = 92488.17999999999   # I99: 101658.3 -> 92488.17999999999
$ws.Cells.Item(99, 11).Value = 92488.17999999999   # K99: 101658.3 -> 92488.17999999999
$ws.Cells.Item(99, 13).Value = -90990.17999999999   # M99: -100160.3 -> -90990.17999999999
$ws.Cells.Item(135, 8).Value = 52849   # H135: 36400 -> 52849
$ws.Cells.Item(135, 10).Value = 52849   # J135: 36400 -> 52849
$ws.Cells.Item(135, 12).Value = 52849   # L135: 36400 -> 52849
$ws.Cells.Item(135, 14).Value = -62989   # N135: -46540 -> -62989
$ws.Cells.Item(11, 14).ClearContents()   # N11: remove (was -14280)

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2594.04   # H31: 2852.889 -> 2594.04
$ws.Cells.Item(31, 9).Value = 1545.5625   # I31: 1671.5385 -> 1545.5625
$ws.Cells.Item(31, 10).Value = 4458   # J31: 5924.4 -> 4458
$ws.Cells.Item(31, 11).Value = 1545.5625   # K31: 1671.5385 -> 1545.5625
$ws.Cells.Item(31, 12).Value = 4458   # L31: 5924.4 -> 4458
$ws.Cells.Item(31, 13).Value = -1250.5625   # M31: -1376.5385 -> -1250.5625
$ws.Cells.Item(31, 14).Value = -5048   # N31: -6514.4 -> -5048
$ws.Cells.Item(34, 8).Value = 2594.04   # H34: 2852.889 -> 2594.04
$ws.Cells.Item(34, 9).Value = 1545.5625   # I34: 1671.5385 -> 1545.5625
$ws.Cells.Item(34, 10).Value = 4458   # J34: 5924.4 -> 4458
$ws.Cells.Item(34, 11).Value = 1545.5625   # K34: 1671.5385 -> 1545.5625
$ws.Cells.Item(34, 12).Value = 4458   # L34: 5924.4 -> 4458
$ws.Cells.Item(34, 13).Value = -1343.5625   # M34: -1469.5385 -> -1343.5625
$ws.Cells.Item(34, 14).Value = -4862   # N34: -6328.4 -> -4862
$ws.Cells.Item(35, 8).Value = 1724.8334   # H35: 2446.75 -> 1724.8334
$ws.Cells.Item(35, 9).Value = 1899.8572   # I35: 3293.5 -> 1899.8572
$ws.Cells.Item(35, 10).Value = 1479.8   # J35: 1600 -> 1479.8
$ws.Cells.Item(35, 11).Value = 1899.8572   # K35: 3293.5 -> 1899.8572
$ws.Cells.Item(35, 12).Value = 1479.8   # L35: 1600 -> 1479.8
$ws.Cells.Item(35, 13).Value = -1605.8572   # M35: -2999.5 -> -1605.8572
$ws.Cells.Item(35, 14).Value = -2067.8   # N35: -2188 -> -2067.8
$ws.Cells.Item(132, 8).Value = 2166655.8   # H132: 2676064.5 -> 2166655.8
$ws.Cells.Item(132, 9).Value = 2527237.2   # I132: 3032456.5 -> 2527237.2
$ws.Cells.Item(132, 10).Value = 3166.6667   # J132: 3125 -> 3166.6667
$ws.Cells.Item(132, 11).Value = 7581711.600000001   # K132: 9097369.5 -> 7581711.600000001
$ws.Cells.Item(132, 12).Value = 9500.000100000001   # L132: 9375 -> 9500.000100000001
$ws.Cells.Item(132, 13).Value = -7579181.600000001   # M132: -9094839.5 -> -7579181.600000001
$ws.Cells.Item(132, 14).Value = -14560.0001   # N132: -14435 -> -14560.0001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 18583456   # H4: 17153966 -> 18583456
$ws.Cells.Item(4, 9).Value = 10000148   # I4: 9091050 -> 10000148
$ws.Cells.Item(4, 11).Value = 30000444   # K4: 27273150 -> 30000444
$ws.Cells.Item(4, 13).Value = -30000332   # M4: -27273038 -> -30000332
$ws.Cells.Item(5, 8).Value = 1079.6   # H5: 1021.4545 -> 1079.6
$ws.Cells.Item(5, 9).Value = 1059.2   # I5: 956 -> 1059.2
$ws.Cells.Item(5, 11).Value = 3177.6   # K5: 2868 -> 3177.6
$ws.Cells.Item(5, 13).Value = -3065.6   # M5: -2756 -> -3065.6
$ws.Cells.Item(121, 8).Value = 845.2273   # H121: 806.0476 -> 845.2273
$ws.Cells.Item(121, 9).Value = 492.15384   # I121: 290 -> 492.15384
$ws.Cells.Item(121, 10).Value = 1355.2222   # J121: 1494.1111 -> 1355.2222
$ws.Cells.Item(121, 11).Value = 1476.46152   # K121: 870 -> 1476.46152
$ws.Cells.Item(121, 12).Value = 4065.6666   # L121: 4482.3333 -> 4065.6666
$ws.Cells.Item(121, 13).Value = -166.4615200000001   # M121: 440 -> -166.4615200000001
$ws.Cells.Item(121, 14).Value = -6685.6666   # N121: -7102.3333 -> -6685.6666
$ws.Cells.Item(132, 8).Value = 4579   # H132: 5982.5 -> 4579
$ws.Cells.Item(132, 9).Value = 2196.75   # I132: 2450 -> 2196.75
$ws.Cells.Item(132, 10).Value = 5940.2856   # J132: 6865.625 -> 5940.2856
$ws.Cells.Item(132, 11).Value = 19770.75   # K132: 22050 -> 19770.75
$ws.Cells.Item(132, 12).Value = 53462.5704   # L132: 61790.625 -> 53462.5704
$ws.Cells.Item(132, 13).Value = -17240.75   # M132: -19520 -> -17240.75
$ws.Cells.Item(132, 14).Value = -58522.5704   # N132: -66850.625 -> -58522.5704
$ws.Cells.Item(135, 8).Value = 1079.6   # H135: 1021.4545 -> 1079.6
$ws.Cells.Item(135, 9).Value = 1059.2   # I135: 956 -> 1059.2
$ws.Cells.Item(135, 11).Value = 9532.800000000001   # K135: 8604 -> 9532.800000000001
$ws.Cells.Item(135, 13).Value = -6997.800000000001   # M135: -6069 -> -6997.800000000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 59490.727   # H70: 110940 -> 59490.727
$ws.Cells.Item(70, 9).Value = 111800   # I70: 137000 -> 111800
$ws.Cells.Item(70, 10).Value = 15899.667   # J70: 6700 -> 15899.667
$ws.Cells.Item(70, 11).Value = 111800   # K70: 137000 -> 111800
$ws.Cells.Item(70, 12).Value = 15899.667   # L70: 6700 -> 15899.667
$ws.Cells.Item(70, 13).Value = -111530   # M70: -136730 -> -111530
$ws.Cells.Item(70, 14).Value = -16439.667   # N70: -7240 -> -16439.667
$ws.Cells.Item(73, 8).Value = 59490.727   # H73: 110940 -> 59490.727
$ws.Cells.Item(73, 9).Value = 111800   # I73: 137000 -> 111800
$ws.Cells.Item(73, 10).Value = 15899.667   # J73: 6700 -> 15899.667
$ws.Cells.Item(73, 11).Value = 111800   # K73: 137000 -> 111800
$ws.Cells.Item(73, 12).Value = 15899.667   # L73: 6700 -> 15899.667
$ws.Cells.Item(73, 13).Value = -110864   # M73: -136064 -> -110864
$ws.Cells.Item(73, 14).Value = -17771.667   # N73: -8572 -> -17771.667
$ws.Cells.Item(102, 8).Value = 1630.75   # H102: 1583.6 -> 1630.75
$ws.Cells.Item(102, 9).Value = 1630.75   # I102: 1583.6 -> 1630.75
$ws.Cells.Item(102, 11).Value = 1630.75   # K102: 1583.6 -> 1630.75
$ws.Cells.Item(102, 13).Value = -8.75   # M102: 38.40000000000009 -> -8.75
$ws.Cells.Item(122, 8).Value = 7340740.5   # H122: 10009146 -> 7340740.5
$ws.Cells.Item(122, 9).Value = 9173968   # I122: 11008461 -> 9173968
$ws.Cells.Item(122, 10).Value = 7833.3335   # J122: 16000 -> 7833.3335
$ws.Cells.Item(122, 11).Value = 27521904   # K122: 33025383 -> 27521904
$ws.Cells.Item(122, 12).Value = 23500.0005   # L122: 48000 -> 23500.0005
$ws.Cells.Item(122, 13).Value = -27519454   # M122: -33022933 -> -27519454
$ws.Cells.Item(122, 14).Value = -28400.0005   # N122: -52900 -> -28400.0005
$ws.Cells.Item(126, 8).Value = 7750   # H126: 3773.0715 -> 7750
$ws.Cells.Item(126, 9).Value = 0   # I126: 2479.8 -> 0
$ws.Cells.Item(126, 10).Value = 7750   # J126: 4491.5557 -> 7750
$ws.Cells.Item(126, 11).Value = 0   # K126: 7439.400000000001 -> 0
$ws.Cells.Item(126, 12).Value = 23250   # L126: 13474.6671 -> 23250
$ws.Cells.Item(126, 14).Value = -28190   # N126: -18414.6671 -> -28190
$ws.Cells.Item(132, 8).Value = 4262.778   # H132: 4290.5557 -> 4262.778
$ws.Cells.Item(132, 10).Value = 6213.5713   # J132: 6320.7144 -> 6213.5713
$ws.Cells.Item(132, 12).Value = 18640.7139   # L132: 18962.1432 -> 18640.7139
$ws.Cells.Item(132, 14).Value = -23700.7139   # N132: -24022.1432 -> -23700.7139
$ws.Cells.Item(126, 13).ClearContents()   # M126: remove (was -4969.400000000001)

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 25929094   # H40: 10374041 -> 25929094
$ws.Cells.Item(40, 9).Value = 6004   # I40: 4205.6 -> 6004
$ws.Cells.Item(40, 11).Value = 6004   # K40: 4205.6 -> 6004
$ws.Cells.Item(40, 13).Value = -5868   # M40: -4069.6 -> -5868
$ws.Cells.Item(132, 8).Value = 2764.4783   # H132: 2686.7917 -> 2764.4783
$ws.Cells.Item(132, 9).Value = 2079.5386   # I132: 1995.2858 -> 2079.5386
$ws.Cells.Item(132, 11).Value = 6238.6158   # K132: 5985.857400000001 -> 6238.6158
$ws.Cells.Item(132, 13).Value = -3708.6158   # M132: -3455.857400000001 -> -3708.6158
$ws.Cells.Item(136, 8).Value = 2377.1765   # H136: 2122.4666 -> 2377.1765
$ws.Cells.Item(136, 9).Value = 1938.7778   # I136: 1778.6666 -> 1938.7778
$ws.Cells.Item(136, 10).Value = 2870.375   # J136: 3497.6667 -> 2870.375
$ws.Cells.Item(136, 11).Value = 5816.3334   # K136: 5335.9998 -> 5816.3334
$ws.Cells.Item(136, 12).Value = 8611.125   # L136: 10493.0001 -> 8611.125
$ws.Cells.Item(136, 13).Value = -3266.3334   # M136: -2785.9998 -> -3266.3334
$ws.Cells.Item(136, 14).Value = -13711.125   # N136: -15593.0001 -> -13711.125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1986.6   # H81: 2004.7 -> 1986.6
$ws.Cells.Item(81, 9).Value = 1045.75   # I81: 1068.375 -> 1045.75
$ws.Cells.Item(81, 11).Value = 2091.5   # K81: 2136.75 -> 2091.5
$ws.Cells.Item(81, 13).Value = -1030.5   # M81: -1075.75 -> -1030.5
$ws.Cells.Item(84, 8).Value = 1986.6   # H84: 2004.7 -> 1986.6
$ws.Cells.Item(84, 9).Value = 1045.75   # I84: 1068.375 -> 1045.75
$ws.Cells.Item(84, 11).Value = 10457.5   # K84: 10683.75 -> 10457.5
$ws.Cells.Item(84, 13).Value = -5153.5   # M84: -5379.75 -> -5153.5
$ws.Cells.Item(107, 8).Value = 1604.0869   # H107: 1604.1305 -> 1604.0869
$ws.Cells.Item(107, 9).Value = 1093.6   # I107: 1093.6666 -> 1093.6
$ws.Cells.Item(107, 11).Value = 3280.8   # K107: 3280.9998 -> 3280.8
$ws.Cells.Item(107, 13).Value = -1360.8   # M107: -1360.9998 -> -1360.8
$ws.Cells.Item(122, 8).Value = 2129.4167   # H122: 2004.0769 -> 2129.4167
$ws.Cells.Item(122, 9).Value = 2226   # I122: 1880.8 -> 2226
$ws.Cells.Item(122, 11).Value = 6678   # K122: 5642.4 -> 6678
$ws.Cells.Item(122, 13).Value = -4228   # M122: -3192.4 -> -4228
$ws.Cells.Item(123, 8).Value = 74799.5   # H123: 74800 -> 74799.5
$ws.Cells.Item(123, 10).Value = 74799.5   # J123: 74800 -> 74799.5
$ws.Cells.Item(123, 12).Value = 74799.5   # L123: 74800 -> 74799.5
$ws.Cells.Item(123, 14).Value = -84599.5   # N123: -84600 -> -84599.5
$ws.Cells.Item(126, 8).Value = 2843.2307   # H126: 2730.2942 -> 2843.2307
$ws.Cells.Item(126, 9).Value = 1771.2   # I126: 1784.5834 -> 1771.2
$ws.Cells.Item(126, 10).Value = 6416.6665   # J126: 5000 -> 6416.6665
$ws.Cells.Item(126, 11).Value = 5313.6   # K126: 5353.7502 -> 5313.6
$ws.Cells.Item(126, 12).Value = 19249.9995   # L126: 15000 -> 19249.9995
$ws.Cells.Item(126, 13).Value = -2843.6   # M126: -2883.7502 -> -2843.6
$ws.Cells.Item(126, 14).Value = -24189.9995   # N126: -19940 -> -24189.9995
$ws.Cells.Item(132, 8).Value = 2706.75   # H132: 2746.5789 -> 2706.75
$ws.Cells.Item(132, 9).Value = 2549   # I132: 2598.9167 -> 2549
$ws.Cells.Item(132, 11).Value = 7647   # K132: 7796.750100000001 -> 7647
$ws.Cells.Item(132, 13).Value = -5117   # M132: -5266.750100000001 -> -5117
